$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 10:17 AM"

# --- Top Gainers sheet: row 8 (GENUSPOWER) value updates ---
$gainers = $wb.Worksheets.Item("Top Gainers")
$gainers.Range("C8").Value = 9.526999999999999
$gainers.Range("D8").Value = 7.7613
$gainers.Range("E8").Value = 4.5638

# --- Top Losers sheet: rows 36-42, 52-57, 76 ---
$losers = $wb.Worksheets.Item("Top Losers")
$losers.Range("B36").Value = "SPARC"
$losers.Range("C36").Value = -3.1709
$losers.Range("D36").Value = 4.8337
$losers.Range("E36").Value = 6.3311

$losers.Range("B37").Value = "NLCINDIA"
$losers.Range("C37").Value = -3.0757
$losers.Range("D37").Value = -4.5618
$losers.Range("E37").Value = -11.6431

$losers.Range("B38").Value = "OAL"
$losers.Range("C38").Value = -3.0452
$losers.Range("D38").Value = -1.3752
$losers.Range("E38").Value = 8.6291

$losers.Range("B39").Value = "MPSLTD"
$losers.Range("C39").Value = -3.0335
$losers.Range("D39").Value = -4.3902
$losers.Range("E39").Value = 2.434

$losers.Range("B40").Value = "DRREDDY"
$losers.Range("C40").Value = -2.9859
$losers.Range("D40").Value = -2.5475
$losers.Range("E40").Value = 2.2228

$losers.Range("B41").Value = "THEMISMED"
$losers.Range("C41").Value = -2.981
$losers.Range("D41").Value = -1.9763
$losers.Range("E41").Value = -8.0187

$losers.Range("B42").Value = "ROSSTECH"
$losers.Range("C42").Value = -2.9778
$losers.Range("D42").Value = 1.9028
$losers.Range("E42").Value = -6.8057

$losers.Range("B52").Value = "PFOCUS"
$losers.Range("C52").Value = -2.7039
$losers.Range("D52").Value = -2.6276
$losers.Range("E52").Value = -1.2163

$losers.Range("B53").Value = "ANANDRATHI"
$losers.Range("C53").Value = -2.6504
$losers.Range("D53").Value = -0.4304
$losers.Range("E53").Value = 9.6646

$losers.Range("B54").Value = "CANHLIFE"
$losers.Range("C54").Value = -2.6148
$losers.Range("D54").Value = "N/A"
$losers.Range("E54").Value = "N/A"

$losers.Range("B55").Value = "GKENERGY"
$losers.Range("C55").Value = -2.6122
$losers.Range("D55").Value = -9.8077
$losers.Range("E55").Value = 23.2758

$losers.Range("B56").Value = "SGFIN"
$losers.Range("C56").Value = -2.592
$losers.Range("D56").Value = -0.0627
$losers.Range("E56").Value = 11.7235

$losers.Range("B57").Value = "ARVINDFASN"
$losers.Range("C57").Value = -2.549
$losers.Range("D57").Value = -2.9892
$losers.Range("E57").Value = -4.4223

$losers.Range("B76").Value = "WEALTH"
$losers.Range("C76").Value = -2.3047
$losers.Range("D76").Value = -3.8606
$losers.Range("E76").Value = -2.8234

# --- 1 Month Performance sheet: rows 10,25,26,30,53-55,58-69,71,72 ---
$perf = $wb.Worksheets.Item("1 Month Performance")
$perf.Range("C10").Value = 51.0181

$perf.Range("B25").Value = "ONMOBILE"
$perf.Range("C25").Value = 34.7681

$perf.Range("B26").Value = "RAMCOSYS"
$perf.Range("C26").Value = 34.6928

$perf.Range("C30").Value = 30.3743

$perf.Range("B53").Value = "PRIVISCL"
$perf.Range("C53").Value = 22.5784

$perf.Range("B54").Value = "CPEDU"
$perf.Range("C54").Value = 22.3786

$perf.Range("B55").Value = "LORDSCHLO"
$perf.Range("C55").Value = 22.1791

$perf.Range("C58").Value = 21.8039

$perf.Range("B59").Value = "ORBTEXP"
$perf.Range("C59").Value = 21.6115

$perf.Range("B60").Value = "GRMOVER"
$perf.Range("C60").Value = 20.2922

$perf.Range("B61").Value = "CEATLTD"
$perf.Range("C61").Value = 20.0239

$perf.Range("B62").Value = "ATL"
$perf.Range("C62").Value = 19.9362

$perf.Range("B63").Value = "SUBROS"
$perf.Range("C63").Value = 19.834

$perf.Range("B64").Value = "HITECHGEAR"
$perf.Range("C64").Value = 19.8096

$perf.Range("B65").Value = "FEDERALBNK"
$perf.Range("C65").Value = 19.6872

$perf.Range("B66").Value = "USHAMART"
$perf.Range("C66").Value = 19.6172

$perf.Range("B67").Value = "BANKINDIA"
$perf.Range("C67").Value = 19.3067

$perf.Range("B68").Value = "REPRO"
$perf.Range("C68").Value = 19.3014

$perf.Range("B69").Value = "RBLBANK"
$perf.Range("C69").Value = 19.2556

$perf.Range("B71").Value = "KARURVYSYA"
$perf.Range("C71").Value = 19.11

$perf.Range("B72").Value = "IIFL"
$perf.Range("C72").Value = 18.9853
